# Scheduled runner update: refresh profit figures across the per-job
# leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Updates recalculated price/profit columns (H,I,J,K,L,M,N) for the
# rows whose underlying market prices changed since the last run.
# A few rows gain a previously-empty HQ profit figure in column N,
# and one row's HQ profit figure is cleared because it is no longer
# applicable.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1761.6666
$ws.Range("I15").Value = 1761.6666
$ws.Range("K15").Value = 5284.9998
$ws.Range("M15").Value = -5115.9998

$ws.Range("H94").Value = 3299.6667
$ws.Range("I94").Value = 3299.6667
$ws.Range("K94").Value = 3299.6667
$ws.Range("M94").Value = -2848.6667

$ws.Range("H129").Value = 135342.4
$ws.Range("I129").Value = 201563.7
$ws.Range("K129").Value = 604691.1000000001
$ws.Range("M129").Value = -599691.1000000001

$ws.Range("H130").Value = 28099.8
$ws.Range("J130").Value = 28099.8
$ws.Range("L130").Value = 28099.8
$ws.Range("N130").Value = -38139.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 23841.857
$ws.Range("I2").Value = 44790.145
$ws.Range("K2").Value = 44790.145
$ws.Range("M2").Value = -44677.145

$ws.Range("H32").Value = 24394786
$ws.Range("I32").Value = 26320414
$ws.Range("K32").Value = 26320414
$ws.Range("M32").Value = -26320127

$ws.Range("H74").Value = 1358.5883
$ws.Range("I74").Value = 1438.3
$ws.Range("J74").Value = 1068.7273
$ws.Range("K74").Value = 1438.3
$ws.Range("L74").Value = 1068.7273
$ws.Range("M74").Value = -564.3
$ws.Range("N74").Value = -2816.7273

$ws.Range("H77").Value = 1358.5883
$ws.Range("I77").Value = 1438.3
$ws.Range("J77").Value = 1068.7273
$ws.Range("K77").Value = 7191.5
$ws.Range("L77").Value = 5343.636500000001
$ws.Range("M77").Value = -2823.5
$ws.Range("N77").Value = -14079.6365

$ws.Range("H109").Value = 49999.816
$ws.Range("J109").Value = 49999.816
$ws.Range("L109").Value = 49999.816
$ws.Range("N109").Value = -52773.816

$ws.Range("H114").Value = 95000
$ws.Range("J114").Value = 95000
$ws.Range("L114").Value = 95000
$ws.Range("N114").Value = -103678

$ws.Range("H116").Value = 23841.857
$ws.Range("I116").Value = 44790.145
$ws.Range("K116").Value = 44790.145
$ws.Range("M116").Value = -42496.145

$ws.Range("H122").Value = 4336.7144
$ws.Range("J122").Value = 5428.091
$ws.Range("L122").Value = 16284.273
$ws.Range("N122").Value = -21184.273

$ws.Range("H132").Value = 2210.32
$ws.Range("I132").Value = 2210.32
$ws.Range("K132").Value = 6630.960000000001
$ws.Range("M132").Value = -4100.960000000001

$ws.Range("H135").Value = 69535.336
$ws.Range("J135").Value = 69535.336
$ws.Range("L135").Value = 69535.336
$ws.Range("N135").Value = -79675.336

$ws.Range("H139").Value = 79545.45
$ws.Range("J139").Value = 79545.45
$ws.Range("L139").Value = 79545.45
$ws.Range("N139").Value = -89825.45

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 23841.857
$ws.Range("I3").Value = 44790.145
$ws.Range("K3").Value = 44790.145
$ws.Range("M3").Value = -44676.145

$ws.Range("H20").Value = 2003.826
$ws.Range("I20").Value = 1321.9333
$ws.Range("J20").Value = 3282.375
$ws.Range("K20").Value = 1321.9333
$ws.Range("L20").Value = 3282.375
$ws.Range("M20").Value = -1074.9333
$ws.Range("N20").Value = -3776.375

$ws.Range("H81").Value = 82328.336
$ws.Range("J81").Value = 82328.336
$ws.Range("L81").Value = 82328.336
$ws.Range("N81").Value = -84450.336

$ws.Range("H84").Value = 82328.336
$ws.Range("J84").Value = 82328.336
$ws.Range("L84").Value = 246985.008
$ws.Range("N84").Value = -257593.008

$ws.Range("H94").Value = 2526.3928
$ws.Range("I94").Value = 2347.7778
$ws.Range("J94").Value = 2847.9
$ws.Range("K94").Value = 2347.7778
$ws.Range("L94").Value = 2847.9
$ws.Range("M94").Value = -1896.7778
$ws.Range("N94").Value = -3749.9

$ws.Range("H134").Value = 2858.1667
$ws.Range("I134").Value = 2340.4375
$ws.Range("K134").Value = 7021.3125
$ws.Range("M134").Value = -4486.3125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 661.625
$ws.Range("I22").Value = 630.6
$ws.Range("K22").Value = 630.6
$ws.Range("M22").Value = -280.6

$ws.Range("H58").Value = 3050.5454
$ws.Range("I58").Value = 2092.6667
$ws.Range("K58").Value = 2092.6667
$ws.Range("M58").Value = -1889.6667

$ws.Range("H99").Value = 17265118
$ws.Range("I99").Value = 3488789.2
$ws.Range("J99").Value = 33337500
$ws.Range("K99").Value = 3488789.2
$ws.Range("L99").Value = 33337500
$ws.Range("M99").Value = -3487291.2
$ws.Range("N99").Value = -33340496

$ws.Range("H122").Value = 491662
$ws.Range("I122").Value = 932372.9399999999
$ws.Range("J122").Value = 6880
$ws.Range("K122").Value = 2797118.82
$ws.Range("L122").Value = 20640
$ws.Range("M122").Value = -2794668.82
$ws.Range("N122").Value = -25540

$ws.Range("H126").Value = 17265118
$ws.Range("I126").Value = 3488789.2
$ws.Range("J126").Value = 33337500
$ws.Range("K126").Value = 10466367.6
$ws.Range("L126").Value = 100012500
$ws.Range("M126").Value = -10463897.6
$ws.Range("N126").Value = -100017440

$ws.Range("H132").Value = 2805.4194
$ws.Range("I132").Value = 2242.96
$ws.Range("J132").Value = 5149
$ws.Range("K132").Value = 6728.88
$ws.Range("L132").Value = 15447
$ws.Range("M132").Value = -4198.88
$ws.Range("N132").Value = -20507

$ws.Range("H134").Value = 7929.222
$ws.Range("I134").Value = 7381.1333
$ws.Range("J134").Value = 10669.667
$ws.Range("K134").Value = 22143.3999
$ws.Range("L134").Value = 32009.001
$ws.Range("M134").Value = -19608.3999
$ws.Range("N134").Value = -37079.001

$ws.Range("H136").Value = 3050.5454
$ws.Range("I136").Value = 2092.6667
$ws.Range("K136").Value = 6278.000100000001
$ws.Range("M136").Value = -3728.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 69503624
$ws.Range("I4").Value = 69503624
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 208510872
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -208510760
$ws.Range("N4").ClearContents()

$ws.Range("H12").Value = 2.9166667

$ws.Range("H132").Value = 3136.182
$ws.Range("I132").Value = 2799.6
$ws.Range("K132").Value = 25196.4
$ws.Range("M132").Value = -22666.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 557305.5
$ws.Range("I70").Value = 1111111
$ws.Range("J70").Value = 3500
$ws.Range("K70").Value = 1111111
$ws.Range("L70").Value = 3500
$ws.Range("M70").Value = -1110841
$ws.Range("N70").Value = -4040

$ws.Range("H73").Value = 557305.5
$ws.Range("I73").Value = 1111111
$ws.Range("J73").Value = 3500
$ws.Range("K73").Value = 1111111
$ws.Range("L73").Value = 3500
$ws.Range("M73").Value = -1110175
$ws.Range("N73").Value = -5372

$ws.Range("H80").Value = 77946.5
$ws.Range("I80").Value = 164669
$ws.Range("K80").Value = 164669
$ws.Range("M80").Value = -163671

$ws.Range("H83").Value = 77946.5
$ws.Range("I83").Value = 164669
$ws.Range("K83").Value = 823345
$ws.Range("M83").Value = -818353

$ws.Range("H102").Value = 3794.5
$ws.Range("I102").Value = 2771.6316
$ws.Range("K102").Value = 2771.6316
$ws.Range("M102").Value = -1149.6316

$ws.Range("H119").Value = 90000
$ws.Range("J119").Value = 90000
$ws.Range("L119").Value = 90000
$ws.Range("N119").Value = -99676

$ws.Range("H122").Value = 5521.5293
$ws.Range("I122").Value = 5126.8887
$ws.Range("J122").Value = 5965.5
$ws.Range("K122").Value = 15380.6661
$ws.Range("L122").Value = 17896.5
$ws.Range("M122").Value = -12930.6661
$ws.Range("N122").Value = -22796.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 722.0833
$ws.Range("I16").Value = 670.75
$ws.Range("K16").Value = 670.75
$ws.Range("M16").Value = -500.75

$ws.Range("H22").Value = 945.95
$ws.Range("I22").Value = 719.8
$ws.Range("J22").Value = 1021.3333
$ws.Range("K22").Value = 719.8
$ws.Range("L22").Value = 1021.3333
$ws.Range("M22").Value = -424.8
$ws.Range("N22").Value = -1611.3333

$ws.Range("H27").Value = 945.95
$ws.Range("I27").Value = 719.8
$ws.Range("J27").Value = 1021.3333
$ws.Range("K27").Value = 719.8
$ws.Range("L27").Value = 1021.3333
$ws.Range("M27").Value = -612.8
$ws.Range("N27").Value = -1235.3333

$ws.Range("H40").Value = 7379.8096
$ws.Range("I40").Value = 8528.857
$ws.Range("K40").Value = 8528.857
$ws.Range("M40").Value = -8392.857

$ws.Range("H43").Value = 535500
$ws.Range("I43").Value = 10000
$ws.Range("J43").Value = 563157.9
$ws.Range("K43").Value = 10000
$ws.Range("L43").Value = 563157.9
$ws.Range("M43").Value = -9807
$ws.Range("N43").Value = -563543.9

$ws.Range("H55").Value = 331.94736
$ws.Range("I55").Value = 485
$ws.Range("J55").Value = 161.88889
$ws.Range("K55").Value = 485
$ws.Range("L55").Value = 161.88889
$ws.Range("M55").Value = -312
$ws.Range("N55").Value = -507.88889

$ws.Range("H132").Value = 2270.2888
$ws.Range("I132").Value = 1865.4117
$ws.Range("J132").Value = 3521.7273
$ws.Range("K132").Value = 5596.2351
$ws.Range("L132").Value = 10565.1819
$ws.Range("M132").Value = -3066.2351
$ws.Range("N132").Value = -15625.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 7303.3335
$ws.Range("I29").Value = 14010
$ws.Range("J29").Value = 3950
$ws.Range("K29").Value = 14010
$ws.Range("L29").Value = 3950
$ws.Range("M29").Value = -13720
$ws.Range("N29").Value = -4530

$ws.Range("H122").Value = 1978.8889
$ws.Range("I122").Value = 1816.0714
$ws.Range("J122").Value = 2548.75
$ws.Range("K122").Value = 5448.2142
$ws.Range("L122").Value = 7646.25
$ws.Range("M122").Value = -2998.2142
$ws.Range("N122").Value = -12546.25

$ws.Range("H126").Value = 2850.5
$ws.Range("I126").Value = 2972
$ws.Range("K126").Value = 8916
$ws.Range("M126").Value = -6446
